$wb = $excel.ActiveWorkbook

# Source sheet to copy header/row-label cell formatting from (bold, bordered, centered style)
$srcSheet = $wb.Worksheets.Item(1)

# Add the new worksheet after the last existing sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "n1000_f_init5_cont0_disc5_sep5p"

# Match the outlinePr settings (summaryBelow / summaryRight) used on the other sheets
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Copy the styled "mean" / "std" header formatting (row1: A1/B1 use the styled format)
$srcSheet.Range("B5:C5").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)

# Copy the styled row-label formatting used for metric names
$srcSheet.Range("A6:A11").Copy()
$ws.Range("A9:A14").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Header row
$ws.Range("B8").Value = "mean"
$ws.Range("C8").Value = "std"

# Metric rows
$ws.Range("A9").Value = "Test Loss (BCE)"
$ws.Range("B9").Value = 0.0155
$ws.Range("C9").Value = 0.0118

$ws.Range("A10").Value = "Accuracy"
$ws.Range("B10").Value = 0.992
$ws.Range("C10").Value = 0.0104

$ws.Range("A11").Value = "F1-Score"
$ws.Range("B11").Value = 0.9921
$ws.Range("C11").Value = 0.0101

$ws.Range("A12").Value = "Precision"
$ws.Range("B12").Value = 0.9865
$ws.Range("C12").Value = 0.0197

$ws.Range("A13").Value = "Recall"
$ws.Range("B13").Value = 0.998
$ws.Range("C13").Value = 0.0045

$ws.Range("A14").Value = "AUC"
$ws.Range("B14").Value = 0.992
$ws.Range("C14").Value = 0.0104

# Match the original "A1 selected" view state
$ws.Range("A1").Select()

# Restore the originally active sheet/tab so the workbook-level view state
# (activeTab) is left exactly as it was before this edit.
$srcSheet.Activate()
$srcSheet.Range("A1").Select()
